# Completed all basic functionality
# Applies the CommonDialogs.xlsx data update: new DialogIDs (600011-600037),
# new DialogTag/Dialog text pairs, widened column C, and updated selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 600011
$ws.Range("A14").Value = 600012
$ws.Range("A15").Value = 600013
$ws.Range("A16").Value = 600014
$ws.Range("A17").Value = 600015
$ws.Range("A18").Value = 600016
$ws.Range("A19").Value = 600017
$ws.Range("A20").Value = 600018
$ws.Range("A21").Value = 600019
$ws.Range("A22").Value = 600020
$ws.Range("A23").Value = 600021
$ws.Range("A24").Value = 600022
$ws.Range("A25").Value = 600023
$ws.Range("B24").Value = 'room-active-node'
$ws.Range("B25").Value = 'room-inactive-node'
$ws.Range("C24").Value = 'There is a {nodename} which is working'
$ws.Range("C25").Value = 'There is a {nodename} which is not working'
$ws.Range("A26").Value = 600024
$ws.Range("B26").Value = 'action-unable-rn'
$ws.Range("C26").Value = 'I am not able to perform that action right now'
$ws.Range("A27").Value = 600025
$ws.Range("B27").Value = 'player-suicide-by-knife'
$ws.Range("C27").Value = 'I have stabbed myself.'
$ws.Range("A28").Value = 600026
$ws.Range("B28").Value = 'game-end'
$ws.Range("C28").Value = 'Game ends.'
$ws.Range("A29").Value = 600027
$ws.Range("B29").Value = 'npc-murdered-by-knife'
$ws.Range("C29").Value = '{npcname} has stabbed himself.'
$ws.Range("A30").Value = 600028
$ws.Range("B30").Value = 'player-escapes-pod'
$ws.Range("C30").Value = 'I have escaped the spacecraft via the escape pod.'
$ws.Range("A31").Value = 600029
$ws.Range("B31").Value = 'player-suicide-by-gun'
$ws.Range("C31").Value = 'I have shot myself.'
$ws.Range("A32").Value = 600030
$ws.Range("B32").Value = 'npc-murdered-by-gun'
$ws.Range("C32").Value = '{npcname} is killed.'
$ws.Range("A33").Value = 600031
$ws.Range("B33").Value = 'user-death-by-knife'
$ws.Range("C33").Value = '{npcname} has stabbed me.'
$ws.Range("A34").Value = 600032
$ws.Range("B34").Value = 'user-death-by-gun'
$ws.Range("C34").Value = '{npcname} has shot me.'
$ws.Range("A35").Value = 600033
$ws.Range("B35").Value = 'npc-death-by-knife'
$ws.Range("C35").Value = '{npcname}  has stabbed himself.'
$ws.Range("A36").Value = 600034
$ws.Range("B36").Value = 'npc-death-by-gun'
$ws.Range("C36").Value = '{npcname}  has shot himself.  '
$ws.Range("A37").Value = 600035
$ws.Range("B37").Value = 'npc-puts-batteries-pod'
$ws.Range("C37").Value = '{npcname} has placed the batteries in escape pod.'
$ws.Range("A38").Value = 600036
$ws.Range("B38").Value = 'npc-escapes-by-pod'
$ws.Range("C38").Value = '{npcname} has escaped the spacecraft via the escape pod.'
$ws.Range("A39").Value = 600037
$ws.Range("B39").Value = 'npc-in-room'
$ws.Range("C39").Value = '{npcname} is here.'

# --- Widen column C to fit the longer dialog text (closest the UI grid allows) ---
$ws.Columns.Item(3).ColumnWidth = 58.666666666666664

# --- Update the view: move the active selection ---
$ws.Range("D9").Select()
